# Auto-generated cell updates derived from the authoritative diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 4.2
$ws.Range("J2").Value = 4.5
$ws.Range("K2").Value = 2.5
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 2.5
$ws.Range("S2").Value = 1.25
$ws.Range("T2").Value = 3.75
$ws.Range("U2").Value = 1.53
$ws.Range("V2").Value = 2.38
$ws.Range("AA2").Value = 29
$ws.Range("AB2").Value = 29
$ws.Range("AC2").Value = 19
$ws.Range("AD2").Value = 8.5
$ws.Range("AH2").Value = 10
$ws.Range("AI2").Value = 10
$ws.Range("AK2").Value = 15
$ws.Range("AM2").Value = 19
$ws.Range("AO2").Value = 21
$ws.Range("AP2").Value = 23
$ws.Range("AR2").Value = 67
$ws.Range("AS2").Value = 126
$ws.Range("AT2").Value = 3.75
$ws.Range("AZ2").Value = 15
$ws.Range("BC2").Value = 81

# Row 3
$ws.Range("I3").Value = 13
$ws.Range("J3").Value = 1.57
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 21
$ws.Range("S3").Value = 1.22
$ws.Range("T3").Value = 4
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 9.5
$ws.Range("Y3").Value = 9.5
$ws.Range("AB3").Value = 26
$ws.Range("AC3").Value = 19
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 301
$ws.Range("AM3").Value = 51
$ws.Range("AN3").Value = 3.4
$ws.Range("AP3").Value = 15
$ws.Range("AT3").Value = 4
$ws.Range("AY3").Value = 41
$ws.Range("BA3").Value = 201
$ws.Range("BC3").Value = 301

# Row 5
$ws.Range("G5").Value = 1.67
$ws.Range("I5").Value = 5.5
$ws.Range("U5").Value = 2.5
$ws.Range("V5").Value = 1.5
$ws.Range("AC5").Value = 7
$ws.Range("AM5").Value = 67

# Row 6
$ws.Range("G6").Value = 2
$ws.Range("I6").Value = 3.7
$ws.Range("J6").Value = 2.63
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11

# Row 7
$ws.Range("Q7").Value = 1.73
$ws.Range("R7").Value = 2.1

# Row 8
$ws.Range("G8").Value = 3.3
$ws.Range("I8").Value = 2.15
$ws.Range("W8").Value = 11
$ws.Range("AN8").Value = 5.5
$ws.Range("AO8").Value = 19

# Row 9
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 3.2
$ws.Range("Q9").Value = 2.1
$ws.Range("R9").Value = 1.73
$ws.Range("BD9").Value = 151

# Row 10
$ws.Range("I10").Value = 4.33
$ws.Range("J10").Value = 2.5
$ws.Range("X10").Value = 8
$ws.Range("Y10").Value = 8.5
$ws.Range("AL10").Value = 41
$ws.Range("BB10").Value = 126

# Row 11
$ws.Range("G11").Value = 1.7
$ws.Range("H11").Value = 3.7
$ws.Range("I11").Value = 4.75
$ws.Range("J11").Value = 2.3
$ws.Range("L11").Value = 5
$ws.Range("Q11").Value = 1.9
$ws.Range("R11").Value = 1.95
$ws.Range("U11").Value = 1.83
$ws.Range("V11").Value = 1.83
$ws.Range("X11").Value = 8
$ws.Range("AA11").Value = 13
$ws.Range("AE11").Value = 17
$ws.Range("AH11").Value = 13
$ws.Range("AL11").Value = 41
$ws.Range("AU11").Value = 8.5
$ws.Range("AY11").Value = 26
$ws.Range("AZ11").Value = 34
$ws.Range("BC11").Value = 251

# Row 13
$ws.Range("G13").Value = 2.22
$ws.Range("I13").Value = 3.2
$ws.Range("J13").Value = 2.82
$ws.Range("K13").Value = 1.98
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 8.38
$ws.Range("S13").Value = 1.42
$ws.Range("T13").Value = 2.45
$ws.Range("U13").Value = 1.72
$ws.Range("V13").Value = 1.9
$ws.Range("W13").Value = 7.4
$ws.Range("X13").Value = 10.75
$ws.Range("Z13").Value = 23
$ws.Range("AA13").Value = 18.5
$ws.Range("AB13").Value = 28
$ws.Range("AC13").Value = 8.5
$ws.Range("AD13").Value = 5.9
$ws.Range("AE13").Value = 13.5
$ws.Range("AG13").Value = 500
$ws.Range("AH13").Value = 9
$ws.Range("AI13").Value = 17
$ws.Range("AJ13").Value = 11.25
$ws.Range("AK13").Value = 45
$ws.Range("AL13").Value = 30
$ws.Range("AM13").Value = 37
$ws.Range("AO13").Value = 11.75
$ws.Range("AP13").Value = 20
$ws.Range("AQ13").Value = 50
$ws.Range("AR13").Value = 80
$ws.Range("AS13").Value = 250
$ws.Range("AT13").Value = 2.4
$ws.Range("AU13").Value = 6.8
$ws.Range("AV13").Value = 60
$ws.Range("AX13").Value = 5.1
$ws.Range("AY13").Value = 18
$ws.Range("AZ13").Value = 25
$ws.Range("BC13").Value = 350
